$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.620731830596924
$ws.Range("B1").Value = 2.749539136886597
$ws.Range("C1").Value = 2.322569608688354
$ws.Range("D1").Value = 2.424772500991821
$ws.Range("E1").Value = 2.717814445495605
